# Translate the English template into Arabic, matching the target diff.
$d = $word.ActiveDocument

function ReplaceAll($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# 1) Top language-picker line: "English" (hyperlink) -> "الإنجليزية"
#    Both occurrences of the bare word "English" become the same Arabic word,
#    so a global replace-all is safe and unambiguous.
ReplaceAll "English" "الإنجليزية"

# 2) " / Portuguese / French / Thai / Vietnamese / Spanish" -> Arabic list
ReplaceAll " / Portuguese / French / Thai / Vietnamese / Spanish" " /البرتغالية/الفرنسية/التايلندية/الفيتنامية/الإسبانية"

# 3) Brief paragraph
ReplaceAll "An email sent to partners who have attended the event. This email will include a photo gallery It will be sent via customer.io" "تم إرسال بريد إلكتروني إلى الشركاء الذين حضروا الحدث. سيتضمن هذا البريد الإلكتروني معرض صور سيتم إرساله عبر customer.io"

# 4) Target audience value
ReplaceAll "Event attendees" "الحاضرون في الحدث"

# 5) "Subject: " label
ReplaceAll "Subject: " "الموضوع: "

# 6) "Thank you for coming to " (subject line lead-in)
ReplaceAll "Thank you for coming to " "شكرًا لقدومك إلى "

# 7) The FIRST "[EVENT NAME]" (in the Subject line, right after the text we just
#    replaced) becomes "[اسم الحدث]"; the second one (later, in the body
#    paragraph "Thank you for attending [EVENT NAME] in [CITY]...") must stay
#    untouched. Scope the Find to a narrow range right after the subject intro.
$scope = $d.Content
$scope.Find.Execute("شكرًا لقدومك إلى ", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0) | Out-Null
$afterSubjectIntro = $d.Range($scope.End, $scope.End + 40)
$afterSubjectIntro.Find.Execute("[EVENT NAME]", $false, $false, $false, $false, $false, `
                                 $true, 1, $false, "[اسم الحدث]", 2) | Out-Null

# 8) Heading "You made our event a success! 🎉"
ReplaceAll "You made our event a success! 🎉" "لقد ساهمت في نجاح هذا الحدث! 🎉"

# 9) "Hi " greeting lead-in
ReplaceAll "Hi " "مرحبًا  "

# 10) The ", " that follows "[PARTNER NAME]" becomes ",، ". "[PARTNER NAME]"
#     itself is unique in the document, so scope the Find right after it.
$scope2 = $d.Content
$scope2.Find.Execute("[PARTNER NAME]", $false, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null
$afterPartnerName = $d.Range($scope2.End, $scope2.End + 10)
$afterPartnerName.Find.Execute(", ", $false, $false, $false, $false, $false, `
                                $true, 1, $false, ",، ", 2) | Out-Null

# 11) Closing sentence of the "Thank you for attending..." paragraph
ReplaceAll ". We hope you had a great time, and it was a pleasure getting to know you!" ". نأمل أن تكونوا قد قضيتوا وقتًا رائعًا، وكان من دواعي سرورنا التعرف عليكم!"

# 12) "We hope the event inspired you..." closing paragraph
ReplaceAll "We hope the event inspired you as much as it did us, and let’s keep growing together!" "نأمل أن يكون هذا الحدث مصدر إلهام لك بقدر ما ألهمنا، ودعونا نستمر في التطور معًا!"

# 13) Comment body text
$c = $d.Comments.Item(1)
$c.Range.Text = "اختر أيًا منهما"

Write-Output "done"
